{"js": "// Office.js (Word JavaScript API) script\n// Applies the \"fix small errors\" edits:\n//   1) \"This module includes 2 micro-modules.\" -> \"This module includes 2 micro modules.\"\n//   2) \"Bell-\" / \"LaPadula\" (spell-checked, split across runs) -> single clean run \"Bell-LaPadula (BLP) Model\"\n//   3) \"Students will be able to know \" -> \"Students will be able to list \" (both occurrences)\n//   4) The \"_GoBack\" bookmark is relocated from the end-of-document page-break paragraph\n//      to right after the word \"list\" in the second \"Students will be able to list\" occurrence.\n\nconst body = context.document.body;\n\n// 1) Fix \"micro-modules\" -> \"micro modules\" (typo/inconsistency fix)\nconst microResults = body.search(\"This module includes 2 micro-modules.\", { matchCase: false });\nmicroResults.load(\"text\");\nawait context.sync();\nif (microResults.items.length > 0) {\n  microResults.items[0].insertText(\"This module includes 2 micro modules.\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Clean up \"Bell-LaPadula (BLP) Model\" (removes stray spell-check run splits / proofErr markers)\nconst blpResults = body.search(\"Bell-LaPadula (BLP) Model\", { matchCase: false });\nblpResults.load(\"text\");\nawait context.sync();\nif (blpResults.items.length > 0) {\n  blpResults.items[0].insertText(\"Bell-LaPadula (BLP) Model\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Replace \"Students will be able to know \" -> \"Students will be able to list \" (both occurrences)\nconst knowResults = body.search(\"Students will be able to know \", { matchCase: false });\nknowResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < knowResults.items.length; i++) {\n  knowResults.items[i].insertText(\"Students will be able to list \", \"Replace\");\n}\nawait context.sync();\n\n// 4) Move the \"_GoBack\" bookmark to sit right after the word \"list\" in the\n//    second \"Students will be able to list\" occurrence (matching where Word\n//    last left the cursor after the edit).\nconst listResults = body.search(\"Students will be able to list\", { matchCase: false });\nlistResults.load(\"text\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_GoBack\");\n\nif (listResults.items.length >= 2) {\n  const secondMatchEnd = listResults.items[1].getRange(\"End\");\n  secondMatchEnd.insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n", "ps1": "# Word COM interop script\n# Applies the \"fix small errors\" edits:\n#   1) \"This module includes 2 micro-modules.\" -> \"This module includes 2 micro modules.\"\n#   2) \"Bell-\" / \"LaPadula\" (spell-checked, split across runs) -> single clean run \"Bell-LaPadula (BLP) Model\"\n#   3) \"Students will be able to know \" -> \"Students will be able to list \" (both occurrences)\n#   4) The \"_GoBack\" bookmark is relocated from the end-of-document page-break paragraph\n#      to right after the word \"list\" in the second \"Students will be able to list\" occurrence.\n\n$d = $word.ActiveDocument\n\n# 1) Fix \"micro-modules\" -> \"micro modules\" (typo/inconsistency fix)\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n[void]$r1.Find.Execute(\"This module includes 2 micro-modules.\", $false, $false, $false, $false, $false, $true, 1, $false, \"This module includes 2 micro modules.\", 2)\n\n# 2) Clean up \"Bell-LaPadula (BLP) Model\" (removes stray spell-check run splits / proofErr markers)\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n[void]$r2.Find.Execute(\"Bell-LaPadula (BLP) Model\", $false, $false, $false, $false, $false, $true, 1, $false, \"Bell-LaPadula (BLP) Model\", 2)\n\n# 3) & 4) Replace \"Students will be able to know \" with \"Students will be able to list \"\n#    one occurrence at a time, and move the \"_GoBack\" bookmark to sit right after the\n#    word \"list\" in the second occurrence (matching where Word last left the cursor).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$newPhrase = \"Students will be able to list\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Students will be able to know \"\n$count = 0\n$bookmarkPos = $null\nwhile ([bool]$rng.Find.Execute()) {\n  $count = $count + 1\n  $rng.Text = \"Students will be able to list \"\n  if ($count -eq 2) {\n    $bookmarkPos = $rng.Start + $newPhrase.Length\n  }\n  $rng.Collapse(0)\n  $rng.End = $d.Content.End\n  if ($count -ge 10) { break }\n}\n\nif ($bookmarkPos -ne $null) {\n  $ins = $d.Range($bookmarkPos, $bookmarkPos)\n  $d.Bookmarks.Add(\"_GoBack\", $ins)\n}\n"}
